{"js": "// Add a new \"Alarm tigger\" section (heading + body paragraph) right after\n// the existing \"Guard Spawner\" section at the end of the document body,\n// preceded by one blank paragraph \u2014 matching the authored OOXML diff.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The last paragraph in the body is the \"Guard Spawner\" description\n// paragraph that ends with \"...at a regular interval.\"\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nconst insertionPoint = lastParagraph.getRange(\"End\");\n\nconst flatOpcXml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  // blank paragraph separating the two sections\n  '<w:p/>' +\n  // \"Alarm tigger\" Heading1 paragraph (note: \"tigger\" is the original\n  // author's typo for \"trigger\", kept as-is with spellcheck markup)\n  '<w:p>' +\n  '<w:pPr><w:pStyle w:val=\"Heading1\"/></w:pPr>' +\n  '<w:r><w:lastRenderedPageBreak/><w:t xml:space=\"preserve\">Alarm </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:t>tigger</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '</w:p>' +\n  // description paragraph\n  '<w:p>' +\n  '<w:r><w:t xml:space=\"preserve\">Attach a alarm trigger component to any object with a collider that the player can touch that will trigger the alarm. When the alarm is triggered, all guard </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:t>spawners</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> in the level will spawn guards.</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\ninsertionPoint.insertOoxml(flatOpcXml, \"After\");\nawait context.sync();\n", "ps1": "# Add a new \"Alarm tigger\" section (heading + body paragraph) right after\n# the existing \"Guard Spawner\" section at the end of the document, preceded\n# by one blank paragraph - matching the authored OOXML diff.\n\n$d = $word.ActiveDocument\n\n# Collapsed range at the very end of the document (after the last\n# paragraph mark) so inserting XML there appends rather than replacing\n# the last paragraph's content.\n$r = $d.Range($d.Content.End, $d.Content.End)\n\n$xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  # blank paragraph separating the two sections\n  '<w:p/>' +\n  # \"Alarm tigger\" Heading1 paragraph (note: \"tigger\" is the original\n  # author's typo for \"trigger\", kept as-is with spellcheck markup)\n  '<w:p>' +\n  '<w:pPr><w:pStyle w:val=\"Heading1\"/></w:pPr>' +\n  '<w:r><w:lastRenderedPageBreak/><w:t xml:space=\"preserve\">Alarm </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:t>tigger</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '</w:p>' +\n  # description paragraph\n  '<w:p>' +\n  '<w:r><w:t xml:space=\"preserve\">Attach a alarm trigger component to any object with a collider that the player can touch that will trigger the alarm. When the alarm is triggered, all guard </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:t>spawners</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> in the level will spawn guards.</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>'\n\n$r.InsertXML($xml)\n"}
